$d = $word.ActiveDocument

# --- 1. "Departement / School of Arts: Wetenschap en Techniek" paragraph ---
# Drop the " / School of Arts" segment and collapse the trailing
# "colon + space" down to a single space, so the heading now reads
# "Departement  Wetenschap en Techniek" (School of Arts removed).
$p4 = $d.Paragraphs.Item(4)
$p4Start = $p4.Range.Start

$schoolOfArts = " / School of Arts"
$rngSchool = $d.Range($p4Start + 11, $p4Start + 11 + $schoolOfArts.Length)
$rngSchool.Text = " "

$rngColon = $d.Range($p4Start + 12, $p4Start + 13)
$rngColon.Text = " "

$rngExtraSpace = $d.Range($p4Start + 13, $p4Start + 14)
$rngExtraSpace.Delete()

# --- 2. "Opleiding: " paragraph: append the study programme name ---
$p5 = $d.Paragraphs.Item(5)
$p5End = $p5.Range.End
$insertPoint = $d.Range($p5End - 1, $p5End - 1)
# Append a harmless placeholder character too, so the insertion point used
# for the bookmark below never sits exactly on the paragraph's last
# character slot (collapsed ranges anchored there are mis-resolved by this
# host's bookmark minting) - it is stripped again right after.
$insertPoint.InsertAfter("elektronica-ictZ")

# --- 3. Relocate the (hidden) _GoBack bookmark to sit right after
#        "elektronica-ict", i.e. at the end of that paragraph - it used to
#        mark the last edit position further down the document. ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$p5 = $d.Paragraphs.Item(5)
$p5End = $p5.Range.End
$bmPos = $d.Range($p5End - 2, $p5End - 2)
$d.Bookmarks.Add("_GoBack", $bmPos)

# Strip the placeholder "Z" again.
$p5 = $d.Paragraphs.Item(5)
$p5End = $p5.Range.End
$placeholder = $d.Range($p5End - 2, $p5End - 1)
$placeholder.Delete()
